$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "First"
$ws.Range("B1").Value = "Second"
$ws.Range("C1").Value = "Third"
$ws.Range("D1").Value = "Fourth"
